# Update "想去人数" (wish-to-go count) figures in the F column of the
# "展览" and "全部类型" sheets to the newly scraped totals.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new value }
$updates = @{
    "展览"   = @{ 2 = 78; 3 = 3930; 4 = 2318; 5 = 458; 7 = 28; 8 = 190; 10 = 5; 11 = 115; 12 = 1458; 14 = 2644 }
    "全部类型" = @{ 2 = 78; 3 = 3930; 4 = 2318; 5 = 458; 7 = 28; 9 = 190; 11 = 5; 12 = 115; 15 = 1458; 17 = 2644 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
